# Commit: "renamed repo, fixed output folder path"
#
# The underlying prediction data didn't change; three extra rows that had
# been duplicated into sheet 1 (rows 8-10) and twelve extra duplicated rows
# in sheet 2 (rows 11-22) are removed, with the remaining rows shifting up
# to close the gap. Excel updates each sheet's <dimension> automatically.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet "s__Schaedlerella glycyrrhizinilytica-b-p": drop the 3 duplicate rows
$ws1.Rows("8:10").Delete()

# Sheet "s__Schaedlerella sp900066545-b-p": drop the 12 duplicate rows
$ws2.Rows("11:22").Delete()
